{"js": "// The author's edit (7/18/2021 \"chuhan CV updated\") makes three\n// visible wording tweaks inside the Work/Project-experience bullets:\n//\n//   1. \"...using python Flask  \"                      -> \"...using Python Flask  \"   (capitalize \"python\")\n//   2. \"...to create awesome user interface...\"        -> \"...to create an awesome user interface...\"\n//   3. \"...including main search page and restaurant...\" -> \"...including the main search page and restaurant...\"\n//\n// (The rest of the diff \u2014 merging runs that used to be split only for\n// <w:proofErr> spell-check markers around \"SKIJoin\", \"useState\", \"Keras\",\n// \"Github\" \u2014 does not change any visible text, and the trailing VML\n// bullet-icon id/size tweak in numbering.xml is a Word-internal save\n// artifact, not a user edit, so neither is reproduced here.)\n\n// 1) \"python Flask\" -> \"Python Flask\"\nlet r1 = context.document.body.search(\"python Flask\", { matchCase: true });\nr1.load(\"text\");\nawait context.sync();\nif (r1.items.length > 0) {\n  r1.items[0].insertText(\"Python Flask\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) \"to create awesome user interface\" -> \"to create an awesome user interface\"\nlet r2 = context.document.body.search(\"to create awesome user interface\", { matchCase: true });\nr2.load(\"text\");\nawait context.sync();\nif (r2.items.length > 0) {\n  r2.items[0].insertText(\"to create an awesome user interface\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) \"including main search page\" -> \"including the main search page\"\nlet r3 = context.document.body.search(\"including main search page\", { matchCase: true });\nr3.load(\"text\");\nawait context.sync();\nif (r3.items.length > 0) {\n  r3.items[0].insertText(\"including the main search page\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The author's edit (7/18/2021 \"chuhan CV updated\") makes three\n# visible wording tweaks inside the Work/Project-experience bullets:\n#\n#   1. \"...using python Flask  \"                        -> \"...using Python Flask  \"   (capitalize \"python\")\n#   2. \"...to create awesome user interface...\"          -> \"...to create an awesome user interface...\"\n#   3. \"...including main search page and restaurant...\" -> \"...including the main search page and restaurant...\"\n#\n# (The rest of the diff only merges runs that used to be split purely for\n# <w:proofErr> spell-check markers around \"SKIJoin\", \"useState\", \"Keras\",\n# \"Github\" -- no visible text changes there -- and the trailing VML\n# bullet-icon id/size tweak in numbering.xml is a Word-internal save\n# artifact, not a user edit, so neither is reproduced here.)\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1) \"python Flask\" -> \"Python Flask\"\nReplace-Text \"python Flask\" \"Python Flask\"\n\n# 2) \"to create awesome user interface\" -> \"to create an awesome user interface\"\nReplace-Text \"to create awesome user interface\" \"to create an awesome user interface\"\n\n# 3) \"including main search page\" -> \"including the main search page\"\nReplace-Text \"including main search page\" \"including the main search page\"\n"}
